$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7291535139083862
$ws.Range("B1").Value = 1.389246940612793
$ws.Range("C1").Value = 4.426088809967041
$ws.Range("D1").Value = 1.831549644470215
$ws.Range("E1").Value = 1.081552624702454
